$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Perejil at Terminal
# Hortofrutícola Agro Chillán. It belongs chronologically right before the
# existing row 49 entry, so insert a fresh row there and push every
# following record (old rows 49-115) down by one (new rows 50-116).
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new observation.
$ws.Range("A49").Value = 7
$ws.Range("B49").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C49").Value = "Ñuble"
$ws.Range("D49").Value = 45174
$ws.Range("E49").Value = 16
$ws.Range("F49").Value = 100112044
$ws.Range("G49").Value = "Perejil"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 200
$ws.Range("K49").Value = 1500
$ws.Range("L49").Value = 1500
$ws.Range("M49").Value = 1500
$ws.Range("N49").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O49").Value = "Región del Maule"
$ws.Range("P49").Value = 1500
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = "Hortaliza"
